# DB Handler changed: better matching. Max 6 respected, Storing match rate
# Update Judge1 (column G) and Judge2 (column H) assignments on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H5").Value = "Yiyang Sun"
$ws.Range("H6").Value = "Era Jain"
$ws.Range("H18").Value = "Era Jain"
$ws.Range("G20").Value = "Elizabeth Carter"
$ws.Range("G28").Value = "Elizabeth Carter"
$ws.Range("G30").Value = "Wanliang Shan"
$ws.Range("H32").Value = "Zhen Ma"
$ws.Range("H33").Value = "Zhen Ma"
$ws.Range("H39").Value = "Yiyang Sun"
$ws.Range("G44").Value = "Yiyang Sun"
$ws.Range("H46").Value = "Yiyang Sun"
$ws.Range("G49").Value = "Elizabeth Carter"
$ws.Range("H51").Value = "Era Jain"
$ws.Range("H52").Value = "Wanliang Shan"
$ws.Range("H55").Value = "Yiyang Sun"
$ws.Range("H58").Value = "Jesse Q. Bond"
$ws.Range("H59").Value = "Ruth Chen"
$ws.Range("H65").Value = "Elizabeth Carter"
$ws.Range("H68").Value = "Zhen Ma"
$ws.Range("G69").Value = "Era Jain"
